$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- MOSFET Chip 2x3 Panel RTM ---
# The BOM now reflects a 2x3 panelized fab run: every line-item quantity
# is multiplied by 6 (the panel population), and the Footprint column is
# simplified from full KiCad library paths down to short descriptive
# names ("0603" for the 0603 passives/LED, "TO-220-3" for the MOSFET).

# Footprint column (C) -- shortened names. The "0603" values look like
# numbers, so a leading apostrophe keeps them stored as text (matching
# the quote-prefixed cells in the real file).
$ws.Range("C2").Value = "'0603"
$ws.Range("C4").Value = "TO-220-3"
$ws.Range("C5").Value = "'0603"
$ws.Range("C6").Value = "'0603"
$ws.Range("C7").Value = "'0603"

# Qty column (G) -- x6 for the 2x3 panel.
$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 12
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 6
$ws.Range("G7").Value = 6

# Column widths tightened to fit the new (shorter) contents.
$ws.Columns.Item(1).ColumnWidth = 9.15
$ws.Columns.Item(7).ColumnWidth = 3.15

# Selection moved onto the MOSFET qty cell.
$ws.Range("G4").Select()
